$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at position 35, pushing existing rows 35..100 down to 36..101
# (this is a log sheet, newest entry goes on top of the "Others" section).
$ws.Rows.Item(35).Insert()

# This sheet always persists every column (A..Y) for every row, even when
# blank, so touch each cell on the new row (re-applying the default style is
# enough to make the engine keep the now-blank cell instead of omitting it)
# to match that pattern before writing the real values into R35/S35.
for ($col = 1; $col -le 25; $col++) {
    $ws.Cells.Item(35, $col).Style = "Normal"
}

# Populate the new row's R/S cells with the new log entry.
$ws.Cells.Item(35, 18).Value = "corporate internet share"
$ws.Cells.Item(35, 19).Value = "2024-09-09 11:14:13"
